# Update market-price derived figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the Diabolos_Profits leve-crafting-profit tables, per refreshed Universalis data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 372153.53
$ws.Range("J17").Value = 456479.3
$ws.Range("L17").Value = 1369437.9
$ws.Range("N17").Value = -1369773.9
$ws.Range("H64").Value = 4034.9412
$ws.Range("I64").Value = 3559.4
$ws.Range("K64").Value = 3559.4
$ws.Range("M64").Value = -3311.4
$ws.Range("H67").Value = 4034.9412
$ws.Range("I67").Value = 3559.4
$ws.Range("K67").Value = 3559.4
$ws.Range("M67").Value = -2701.4
$ws.Range("H75").Value = 59999
$ws.Range("I75").Value = 59999
$ws.Range("K75").Value = 59999
$ws.Range("M75").Value = -59063
$ws.Range("H78").Value = 59999
$ws.Range("I78").Value = 59999
$ws.Range("K78").Value = 179997
$ws.Range("M78").Value = -175317
$ws.Range("H113").Value = 66669868
$ws.Range("I113").Value = 333335000
$ws.Range("J113").Value = 3583.1667
$ws.Range("K113").Value = 333335000
$ws.Range("L113").Value = 3583.1667
$ws.Range("M113").Value = -333331746
$ws.Range("N113").Value = -10091.1667
$ws.Range("H132").Value = 13701079
$ws.Range("J132").Value = 8033.2
$ws.Range("L132").Value = 24099.6
$ws.Range("N132").Value = -29159.6
$ws.Range("H133").Value = 99999
$ws.Range("I133").Value = 99999
$ws.Range("K133").Value = 99999
$ws.Range("M133").Value = -94939
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 748.0333000000001
$ws.Range("I2").Value = 748.0333000000001
$ws.Range("K2").Value = 748.0333000000001
$ws.Range("M2").Value = -635.0333000000001
$ws.Range("H32").Value = 12755.39
$ws.Range("I32").Value = 7108.8335
$ws.Range("K32").Value = 7108.8335
$ws.Range("M32").Value = -6821.8335
$ws.Range("H45").Value = 1853175
$ws.Range("I45").Value = 5555555
$ws.Range("K45").Value = 5555555
$ws.Range("M45").Value = -5555178
$ws.Range("H61").Value = 2730.3845
$ws.Range("I61").Value = 1673.2106
$ws.Range("J61").Value = 5599.857
$ws.Range("K61").Value = 1673.2106
$ws.Range("L61").Value = 5599.857
$ws.Range("M61").Value = -1461.2106
$ws.Range("N61").Value = -6023.857
$ws.Range("H110").Value = 25001630
$ws.Range("I110").Value = 29413160
$ws.Range("K110").Value = 29413160
$ws.Range("M110").Value = -29411115
$ws.Range("H116").Value = 748.0333000000001
$ws.Range("I116").Value = 748.0333000000001
$ws.Range("K116").Value = 748.0333000000001
$ws.Range("M116").Value = 1545.9667
$ws.Range("H122").Value = 2529.5
$ws.Range("I122").Value = 2153.6128
$ws.Range("K122").Value = 6460.8384
$ws.Range("M122").Value = -4010.8384
$ws.Range("H136").Value = 2730.3845
$ws.Range("I136").Value = 1673.2106
$ws.Range("J136").Value = 5599.857
$ws.Range("K136").Value = 5019.6318
$ws.Range("L136").Value = 16799.571
$ws.Range("M136").Value = -2469.6318
$ws.Range("N136").Value = -21899.571

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 748.0333000000001
$ws.Range("I3").Value = 748.0333000000001
$ws.Range("K3").Value = 748.0333000000001
$ws.Range("M3").Value = -634.0333000000001
$ws.Range("H9").Value = 15000
$ws.Range("J9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("N9").Value = -15336
$ws.Range("H20").Value = 40941
$ws.Range("I20").Value = 50751.5
$ws.Range("K20").Value = 50751.5
$ws.Range("M20").Value = -50504.5
$ws.Range("H86").Value = 3413.8333
$ws.Range("I86").Value = 2297.6
$ws.Range("K86").Value = 2297.6
$ws.Range("M86").Value = -1174.6
$ws.Range("H89").Value = 3413.8333
$ws.Range("I89").Value = 2297.6
$ws.Range("K89").Value = 11488
$ws.Range("M89").Value = -5872
$ws.Range("H107").Value = 45461556
$ws.Range("I107").Value = 3018.8333
$ws.Range("J107").Value = 100011800
$ws.Range("K107").Value = 3018.8333
$ws.Range("L107").Value = 100011800
$ws.Range("M107").Value = -1098.8333
$ws.Range("N107").Value = -100015640
$ws.Range("H134").Value = 3047.6453
$ws.Range("I134").Value = 2374.6086
$ws.Range("J134").Value = 4982.625
$ws.Range("K134").Value = 7123.825800000001
$ws.Range("L134").Value = 14947.875
$ws.Range("M134").Value = -4588.825800000001
$ws.Range("N134").Value = -20017.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 34999
$ws.Range("I51").Value = 34999
$ws.Range("K51").Value = 34999
$ws.Range("M51").Value = -34263
$ws.Range("H61").Value = 34999
$ws.Range("I61").Value = 34999
$ws.Range("K61").Value = 34999
$ws.Range("M61").Value = -34651
$ws.Range("H62").Value = 44497.25
$ws.Range("I62").Value = 2995.75
$ws.Range("J62").Value = 65248
$ws.Range("K62").Value = 2995.75
$ws.Range("L62").Value = 65248
$ws.Range("M62").Value = -2371.75
$ws.Range("N62").Value = -66496
$ws.Range("H65").Value = 44497.25
$ws.Range("I65").Value = 2995.75
$ws.Range("J65").Value = 65248
$ws.Range("K65").Value = 14978.75
$ws.Range("L65").Value = 326240
$ws.Range("M65").Value = -11858.75
$ws.Range("N65").Value = -332480
$ws.Range("H68").Value = 72250
$ws.Range("I68").Value = 60000
$ws.Range("J68").Value = 84500
$ws.Range("K68").Value = 60000
$ws.Range("L68").Value = 84500
$ws.Range("M68").Value = -59251
$ws.Range("N68").Value = -85998
$ws.Range("H71").Value = 72250
$ws.Range("I71").Value = 60000
$ws.Range("J71").Value = 84500
$ws.Range("K71").Value = 180000
$ws.Range("L71").Value = 253500
$ws.Range("M71").Value = -176256
$ws.Range("N71").Value = -260988
$ws.Range("H107").Value = 520.64
$ws.Range("I107").Value = 486.65216
$ws.Range("K107").Value = 486.65216
$ws.Range("M107").Value = 1433.34784

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H97").Value = 888.6667
$ws.Range("I97").Value = 622.5
$ws.Range("J97").Value = 1066.1111
$ws.Range("K97").Value = 1867.5
$ws.Range("L97").Value = 3198.3333
$ws.Range("M97").Value = -1371.5
$ws.Range("N97").Value = -4190.3333
$ws.Range("H131").Value = 7127.8203
$ws.Range("I131").Value = 1666.5
$ws.Range("K131").Value = 4999.5
$ws.Range("M131").Value = 40.5
$ws.Range("H137").Value = 2899.68
$ws.Range("J137").Value = 3328.2942
$ws.Range("L137").Value = 9984.882599999999
$ws.Range("N137").Value = -20184.8826

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10754.6
$ws.Range("I70").Value = 9999
$ws.Range("J70").Value = 10943.5
$ws.Range("K70").Value = 9999
$ws.Range("L70").Value = 10943.5
$ws.Range("M70").Value = -9729
$ws.Range("N70").Value = -11483.5
$ws.Range("H73").Value = 10754.6
$ws.Range("I73").Value = 9999
$ws.Range("J73").Value = 10943.5
$ws.Range("K73").Value = 9999
$ws.Range("L73").Value = 10943.5
$ws.Range("M73").Value = -9063
$ws.Range("N73").Value = -12815.5
$ws.Range("H122").Value = 181309.84
$ws.Range("I122").Value = 224111.92
$ws.Range("K122").Value = 672335.76
$ws.Range("M122").Value = -669885.76

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 887.8946999999999
$ws.Range("I16").Value = 634.0769
$ws.Range("J16").Value = 1437.8334
$ws.Range("K16").Value = 634.0769
$ws.Range("L16").Value = 1437.8334
$ws.Range("M16").Value = -464.0769
$ws.Range("N16").Value = -1777.8334
$ws.Range("H22").Value = 714.1667
$ws.Range("I22").Value = 295
$ws.Range("J22").Value = 798
$ws.Range("K22").Value = 295
$ws.Range("L22").Value = 798
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = -1388
$ws.Range("H27").Value = 714.1667
$ws.Range("I27").Value = 295
$ws.Range("J27").Value = 798
$ws.Range("K27").Value = 295
$ws.Range("L27").Value = 798
$ws.Range("M27").Value = -188
$ws.Range("N27").Value = -1012
$ws.Range("H40").Value = 2503619.8
$ws.Range("J40").Value = 4608.2
$ws.Range("L40").Value = 4608.2
$ws.Range("N40").Value = -4880.2
$ws.Range("H68").Value = 7174.143
$ws.Range("I68").Value = 4254.75
$ws.Range("K68").Value = 4254.75
$ws.Range("M68").Value = -3505.75
$ws.Range("H71").Value = 7174.143
$ws.Range("I71").Value = 4254.75
$ws.Range("K71").Value = 21273.75
$ws.Range("M71").Value = -17529.75
$ws.Range("H82").Value = 2365.5
$ws.Range("I82").Value = 2197.7144
$ws.Range("K82").Value = 2197.7144
$ws.Range("M82").Value = -1836.7144
$ws.Range("H85").Value = 2365.5
$ws.Range("I85").Value = 2197.7144
$ws.Range("K85").Value = 2197.7144
$ws.Range("M85").Value = -949.7143999999998
$ws.Range("H122").Value = 4809.263
$ws.Range("I122").Value = 5251.636
$ws.Range("K122").Value = 15754.908
$ws.Range("M122").Value = -13304.908
$ws.Range("H136").Value = 4129.224
$ws.Range("I136").Value = 3900.9019
$ws.Range("K136").Value = 11702.7057
$ws.Range("M136").Value = -9152.705699999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 646.9286
$ws.Range("J107").Value = 501.83334
$ws.Range("L107").Value = 1505.50002
$ws.Range("N107").Value = -5345.500019999999
$ws.Range("H136").Value = 4078.577
$ws.Range("I136").Value = 4178.722
$ws.Range("J136").Value = 3853.25
$ws.Range("K136").Value = 12536.166
$ws.Range("L136").Value = 11559.75
$ws.Range("M136").Value = -9986.165999999999
$ws.Range("N136").Value = -16659.75

